# Generate Report for Handoff
#
# The files 16ae8849-..., 2420e8a0-..., 340b87f2-..., 52c24b06-... (rows 4-7
# in both the "zh-cn" and "de-de" localization-status sheets) have finished
# a fresh handoff cycle: their Priority moves from "low" to "ht", and their
# "Latest Handoff Datetime" is bumped forward by 15 seconds to reflect the
# newly generated handoff xliff timestamp.

function Add-FifteenSeconds([string]$timestamp) {
    $datePart = $timestamp.Substring(0, 10)
    $timePart = $timestamp.Substring(11, 8)
    $parts = $timePart.Split(":")
    $h = [int]$parts[0]
    $m = [int]$parts[1]
    $s = [int]$parts[2]

    $s = $s + 15
    if ($s -ge 60) {
        $s = $s - 60
        $m = $m + 1
    }
    if ($m -ge 60) {
        $m = $m - 60
        $h = $h + 1
    }

    $hh = $h.ToString("00")
    $mm = $m.ToString("00")
    $ss = $s.ToString("00")
    return "$datePart $hh`:$mm`:$ss"
}

$wb = $excel.ActiveWorkbook

$sheetNames = @("zh-cn", "de-de")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    for ($row = 4; $row -le 7; $row++) {

        # Priority: low -> ht
        $ws.Range("E" + $row).Value = "ht"

        # Latest Handoff Datetime: +15 seconds
        $current = $ws.Range("H" + $row).Value()
        $newStamp = Add-FifteenSeconds $current
        $ws.Range("H" + $row).Value = $newStamp
    }
}
